$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (rows 2-3) down by one row to make room for the new
# "home" row, without using Rows.Insert() so cell formatting (style) is not
# copied down from the header row.

# Move "tarazzo" (row 3) -> row 4, with value reset to 0
$ws.Range("A4").Value = $ws.Range("A3").Value2
$ws.Range("B4").Value = 0

# Move "tiles" (row 2) -> row 3, with value reset to 0
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = 0

# New row 2: "home" with value 0
$ws.Range("A2").Value = "home"
$ws.Range("B2").Value = 0
